$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("数组")

# Row 11: LeetCode 283 - Move Zeroes
$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = 283
$ws.Cells.Item(11, 3).Value = "给定一个数组 nums，编写一个函数将所有 0 移动到数组的末尾，同时保持非零元素的相对顺序。 `n示例:`n输入: [0,1,0,3,12]`n输出: [1,3,12,0,0] `n 说明:`n 必须在原数组上操作，不能拷贝额外的数组。`n 尽量减少操作次数。 `n Related Topics 数组 双指针"
$ws.Cells.Item(11, 4).Value = "1 定义index字段，只想非0的元素`n2 如果i位置当前元素是0，则迭代继续`n3 如果i位置当前不是0，就将这个元素，移动到index的位置，并且index++`n4 迭代完成数组元素，从index位置到数组终点，填充为0"
$ws.Cells.Item(11, 5).Value = "双指针`n数组覆盖"
$ws.Cells.Item(11, 6).Value = "O(N)"
$ws.Cells.Item(11, 7).Value = "O(1)"

$ws.Rows.Item(11).RowHeight = 198

$ws.Range("D11").Select()
